# Update Haba (Hortaliza) weekly price records: Fecha, Volumen, Precio min/max/prom,
# Origen (Provincia) and Precio $/Kg for rows 2-43 (data rows; header is row 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Row, D(Fecha serial), J(Volumen), K(Precio minimo), L(Precio maximo), M(Precio promedio ponderado), O(Origen), P(Precio $/Kg)
$rows = @(
    @(2, 44484, 400, 9000, 10000, 9500, "Provincia del Elquí", 380),
    @(3, 44370, 520, 13000, 14000, 13500, "Provincia del Elquí", 540),
    @(4, 44473, 500, 8500, 9000, 8750, "Provincia del Elquí", 350),
    @(5, 44880, 560, 7000, 8000, 7500, "Provincia del Elquí", 300),
    @(6, 44858, 500, 7000, 8000, 7500, "Provincia del Elquí", 300),
    @(7, 44799, 500, 10000, 11000, 10500, "Provincia del Elquí", 420),
    @(8, 44797, 1000, 11000, 12000, 11500, "Provincia del Elquí", 460),
    @(9, 44466, 400, 9500, 10000, 9750, "Provincia del Elquí", 390),
    @(10, 44846, 600, 7500, 8000, 7750, "Provincia del Elquí", 310),
    @(11, 44809, 520, 9500, 10000, 9750, "Provincia del Elquí", 390),
    @(12, 44811, 400, 10000, 10500, 10250, "Provincia del Elquí", 410),
    @(13, 44798, 400, 10500, 11000, 10750, "Provincia del Elquí", 430),
    @(14, 44874, 500, 6000, 7000, 6500, "Provincia del Elquí", 260),
    @(15, 44825, 440, 8000, 9000, 8500, "Provincia del Elquí", 340),
    @(16, 44824, 500, 8000, 9000, 8500, "Provincia del Elquí", 340),
    @(17, 44876, 460, 6000, 7000, 6500, "Provincia del Elquí", 260),
    @(18, 44855, 540, 7000, 8000, 7500, "Provincia del Elquí", 300),
    @(19, 44817, 440, 9000, 10000, 9500, "Provincia del Elquí", 380),
    @(20, 44827, 700, 8000, 9000, 8500, "Provincia del Elquí", 340),
    @(21, 44386, 500, 11000, 12000, 11500, "Provincia del Elquí", 460),
    @(22, 44873, 540, 6000, 7000, 6500, "Provincia del Elquí", 260),
    @(23, 44837, 520, 8000, 9000, 8500, "Provincia del Elquí", 340),
    @(24, 44816, 600, 9500, 10000, 9750, "Provincia del Elquí", 390),
    @(25, 44883, 380, 7000, 8000, 7500, "Provincia del Elquí", 300),
    @(26, 44356, 500, 13000, 14000, 13500, "Provincia de Limarí", 540),
    @(27, 44372, 500, 13000, 14000, 13500, "Provincia del Elquí", 540),
    @(28, 44694, 480, 17500, 18000, 17750, "Provincia del Elquí", 710),
    @(29, 44446, 500, 11000, 12000, 11500, "Provincia del Elquí", 460),
    @(30, 44714, 400, 14000, 15000, 14500, "Provincia de Limarí", 580),
    @(31, 44781, 400, 10000, 11000, 10500, "Provincia del Elquí", 420),
    @(32, 44847, 520, 7000, 8000, 7500, "Provincia del Elquí", 300),
    @(33, 44377, 520, 12500, 13000, 12750, "Provincia del Elquí", 510),
    @(34, 44425, 400, 11500, 12000, 11750, "Provincia del Elquí", 470),
    @(35, 44803, 600, 9500, 10000, 9750, "Provincia del Elquí", 390),
    @(36, 44690, 400, 17000, 18000, 17500, "Provincia del Elquí", 700),
    @(37, 44721, 500, 14500, 15000, 14750, "Provincia de Limarí", 590),
    @(38, 44881, 500, 7000, 8000, 7500, "Provincia del Elquí", 300),
    @(39, 44316, 300, 16000, 17000, 16500, "Provincia del Elquí", 660),
    @(40, 44384, 560, 11500, 12000, 11750, "Provincia del Elquí", 470),
    @(41, 44756, 400, 14000, 15000, 14500, "Provincia del Elquí", 580),
    @(42, 44848, 800, 7000, 8000, 7500, "Provincia del Elquí", 300),
    @(43, 44376, 400, 12000, 13000, 12500, "Provincia del Elquí", 500)
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Range("D$row").Value = $r[1]
    $ws.Range("J$row").Value = $r[2]
    $ws.Range("K$row").Value = $r[3]
    $ws.Range("L$row").Value = $r[4]
    $ws.Range("M$row").Value = $r[5]
    $ws.Range("O$row").Value = $r[6]
    $ws.Range("P$row").Value = $r[7]
}
